$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F1").Value = "DateNaissance"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("F2").Value = "51/13/2001"

$ws.Range("F3").Value = 40129
$ws.Range("F3").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("F3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
